# Iteration 3 commit and push
# Adds a new worksheet "Лист3" with a TimeLog table, and updates
# tab-selection / active-cell state on the workbook and on "Лист2".

$wb = $excel.ActiveWorkbook

# --- Add the new worksheet "Лист3" after "Лист2" -----------------------
$sheet2 = $wb.Worksheets.Item("Лист2")
$ws3 = $wb.Worksheets.Add()
$ws3.Name = "Лист3"
$ws3.Move($null, $sheet2)
# Moving invalidates the captured handle in this runtime, so re-fetch it.
$ws3 = $wb.Worksheets.Item("Лист3")

# --- Header row (row 1) --------------------------------------------------
$ws3.Range("A1").Value = "Activity"
$ws3.Range("B1").Value = "Estimated Time"
$ws3.Range("C1").Value = "Real Time"
$ws3.Range("A1:C1").Style = "Accent1"

# --- Data rows -------------------------------------------------------------
$ws3.Range("A3").Value = "Implement The Game / Refactor"
$ws3.Range("B3").Value = "6 hours"
$ws3.Range("C3").Value = "8 hours"

$ws3.Range("A5").Value = "Study the book chapters 8"
$ws3.Range("B5").Value = "2.5 hours "
$ws3.Range("C5").Value = "3 hours"

$ws3.Range("A7").Value = "Create Test Plan"
$ws3.Range("B7").Value = "2 hours"
$ws3.Range("C7").Value = "1.5 hours"

$ws3.Range("A9").Value = "Manual Test"
$ws3.Range("B9").Value = "4 hours"
$ws3.Range("C9").Value = "5 hours"

$ws3.Range("A11").Value = "JUnit5 Auto test"
$ws3.Range("B11").Value = "2 hours"
$ws3.Range("C11").Value = "9 hours"

# --- Column widths (approximate Лист2 layout) -----------------------------
$ws3.Columns.Item(1).ColumnWidth = 27.41
$ws3.Columns.Item(2).ColumnWidth = 14.08

# --- Selection / active sheet state ----------------------------------------
# Лист2 is no longer the active tab; its selection resets to A1.
$sheet2.Range("A1").Select()

# Лист3 becomes the active (selected) sheet with E22 selected.
$ws3.Activate()
$ws3.Range("E22").Select()
